$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header formatting (bold, centered, bordered) from an existing header cell
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply header text values (PasteSpecial for formats only should not
# have touched the text, but set again defensively)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Boolean outlier flags for rows 2-19 (TRUE only for rows 5 and 19)
$outlierRows = @(5, 19)
for ($r = 2; $r -le 19; $r++) {
    $isOutlier = $outlierRows -contains $r
    $ws.Cells.Item($r, 6).Value = $isOutlier
    $ws.Cells.Item($r, 7).Value = $isOutlier
    $ws.Cells.Item($r, 8).Value = $isOutlier
}
